$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 37 and 63: set resultado/profit for newly settled bets ---
$ws.Range("G37").Value = "Fallo"
$ws.Range("H37").Value = -1
$ws.Range("G63").Value = "Fallo"
$ws.Range("H63").Value = -1

# --- Append new match rows 97-118 ---
# Row 97
$ws.Range("A97").Value = 14495004
$ws.Range("B97").NumberFormat = "@"
$ws.Range("B97").Value = "2025-09-03"
$ws.Range("B97").ClearFormats()
$ws.Range("C97").Value = "Karolina Muchova"
$ws.Range("D97").Value = "Naomi Osaka"
$ws.Range("E97").Value = "Gana Karolina Muchova"
$ws.Range("F97").Value = 3

# Row 98
$ws.Range("A98").Value = 14552578
$ws.Range("B98").NumberFormat = "@"
$ws.Range("B98").Value = "2025-09-04"
$ws.Range("B98").ClearFormats()
$ws.Range("C98").Value = "Stefano Travaglia"
$ws.Range("D98").Value = "Gianluca Cadenasso"
$ws.Range("E98").Value = "Gana Gianluca Cadenasso"
$ws.Range("F98").Value = 3.25

# Row 99
$ws.Range("A99").Value = 14552528
$ws.Range("B99").NumberFormat = "@"
$ws.Range("B99").Value = "2025-09-04"
$ws.Range("B99").ClearFormats()
$ws.Range("C99").Value = "Kasidit Samrej"
$ws.Range("D99").Value = "Yu Hsiou Hsu"
$ws.Range("E99").Value = "Gana Kasidit Samrej"
$ws.Range("F99").Value = 4

# Row 100
$ws.Range("A100").Value = 14552614
$ws.Range("B100").NumberFormat = "@"
$ws.Range("B100").Value = "2025-09-04"
$ws.Range("B100").ClearFormats()
$ws.Range("C100").Value = "Marko Topo"
$ws.Range("D100").Value = "Andrew Paulson"
$ws.Range("E100").Value = "Gana Andrew Paulson"
$ws.Range("F100").Value = 2.75

# Row 101
$ws.Range("A101").Value = 14560139
$ws.Range("B101").NumberFormat = "@"
$ws.Range("B101").Value = "2025-09-03"
$ws.Range("B101").ClearFormats()
$ws.Range("C101").Value = "Alina Charaeva"
$ws.Range("D101").Value = "Zhibek Kulambayeva"
$ws.Range("E101").Value = "Gana Zhibek Kulambayeva"
$ws.Range("F101").Value = 3.4

# Row 102
$ws.Range("A102").Value = 14560140
$ws.Range("B102").NumberFormat = "@"
$ws.Range("B102").Value = "2025-09-03"
$ws.Range("B102").ClearFormats()
$ws.Range("C102").Value = "Yexin MA"
$ws.Range("D102").Value = "Veronika Erjavec"
$ws.Range("E102").Value = "Gana Yexin MA"
$ws.Range("F102").Value = 4.5

# Row 103
$ws.Range("A103").Value = 14559638
$ws.Range("B103").NumberFormat = "@"
$ws.Range("B103").Value = "2025-09-03"
$ws.Range("B103").ClearFormats()
$ws.Range("C103").Value = "Maria Kozyreva"
$ws.Range("D103").Value = "Martina Trevisan"
$ws.Range("E103").Value = "Gana Maria Kozyreva"
$ws.Range("F103").Value = 2.75

# Row 104
$ws.Range("A104").Value = 14559642
$ws.Range("B104").NumberFormat = "@"
$ws.Range("B104").Value = "2025-09-03"
$ws.Range("B104").ClearFormats()
$ws.Range("C104").Value = "Ana Sofia Sanchez"
$ws.Range("D104").Value = "Kayla Day"
$ws.Range("E104").Value = "Gana Ana Sofia Sanchez"
$ws.Range("F104").Value = 2.75

# Row 105
$ws.Range("A105").Value = 14559649
$ws.Range("B105").NumberFormat = "@"
$ws.Range("B105").Value = "2025-09-03"
$ws.Range("B105").ClearFormats()
$ws.Range("C105").Value = "Varvara Lepchenko"
$ws.Range("D105").Value = "Alexandra Eala"
$ws.Range("E105").Value = "Gana Varvara Lepchenko"
$ws.Range("F105").Value = 4

# Row 106
$ws.Range("A106").Value = 14559663
$ws.Range("B106").NumberFormat = "@"
$ws.Range("B106").Value = "2025-09-03"
$ws.Range("B106").ClearFormats()
$ws.Range("C106").Value = "Darja Semenistaja"
$ws.Range("D106").Value = "Kaja Juvan"
$ws.Range("E106").Value = "Gana Darja Semenistaja"
$ws.Range("F106").Value = 2.62

# Row 107
$ws.Range("A107").Value = 14559668
$ws.Range("B107").NumberFormat = "@"
$ws.Range("B107").Value = "2025-09-03"
$ws.Range("B107").ClearFormats()
$ws.Range("C107").Value = "Dominika Salkova"
$ws.Range("D107").Value = "Andrea Lazaro Garcia"
$ws.Range("E107").Value = "Gana Andrea Lazaro Garcia"
$ws.Range("F107").Value = 3.25

# Row 108
$ws.Range("A108").Value = 14559667
$ws.Range("B108").NumberFormat = "@"
$ws.Range("B108").Value = "2025-09-03"
$ws.Range("B108").ClearFormats()
$ws.Range("C108").Value = "Tara Würth"
$ws.Range("D108").Value = "Julia Grabher"
$ws.Range("E108").Value = "Gana Tara Würth"
$ws.Range("F108").Value = 3

# Row 109
$ws.Range("A109").Value = 14559660
$ws.Range("B109").NumberFormat = "@"
$ws.Range("B109").Value = "2025-09-04"
$ws.Range("B109").ClearFormats()
$ws.Range("C109").Value = "Darya Astakhova"
$ws.Range("D109").Value = "Anca Todoni"
$ws.Range("E109").Value = "Gana Darya Astakhova"
$ws.Range("F109").Value = 4

# Row 110
$ws.Range("A110").Value = 14600110
$ws.Range("B110").NumberFormat = "@"
$ws.Range("B110").Value = "2025-09-03"
$ws.Range("B110").ClearFormats()
$ws.Range("C110").Value = "Emanuele Mazzeschi"
$ws.Range("D110").Value = "Lorenzo Beraldo"
$ws.Range("E110").Value = "Gana Emanuele Mazzeschi"
$ws.Range("F110").Value = 4.33

# Row 111
$ws.Range("A111").Value = 14600102
$ws.Range("B111").NumberFormat = "@"
$ws.Range("B111").Value = "2025-09-03"
$ws.Range("B111").ClearFormats()
$ws.Range("C111").Value = "Matteo Gribaldo"
$ws.Range("D111").Value = "Juan Cruz Martin Manzano"
$ws.Range("E111").Value = "Gana Matteo Gribaldo"
$ws.Range("F111").Value = 5.5

# Row 112
$ws.Range("A112").Value = 14600109
$ws.Range("B112").NumberFormat = "@"
$ws.Range("B112").Value = "2025-09-03"
$ws.Range("B112").ClearFormats()
$ws.Range("C112").Value = "Raffaele Ciurnelli"
$ws.Range("D112").Value = "Giorgio Tabacco"
$ws.Range("E112").Value = "Gana Raffaele Ciurnelli"
$ws.Range("F112").Value = 4.5

# Row 113
$ws.Range("A113").Value = 14600104
$ws.Range("B113").NumberFormat = "@"
$ws.Range("B113").Value = "2025-09-03"
$ws.Range("B113").ClearFormats()
$ws.Range("C113").Value = "Alessandro Coccioli"
$ws.Range("D113").Value = "Iannis Miletich"
$ws.Range("E113").Value = "Gana Alessandro Coccioli"
$ws.Range("F113").Value = 3.4

# Row 114
$ws.Range("A114").Value = 14600101
$ws.Range("B114").NumberFormat = "@"
$ws.Range("B114").Value = "2025-09-03"
$ws.Range("B114").ClearFormats()
$ws.Range("C114").Value = "Lorenzo Comino"
$ws.Range("D114").Value = "Alberto Morolli"
$ws.Range("E114").Value = "Gana Alberto Morolli"
$ws.Range("F114").Value = 2.5

# Row 115
$ws.Range("A115").Value = 14600108
$ws.Range("B115").NumberFormat = "@"
$ws.Range("B115").Value = "2025-09-03"
$ws.Range("B115").ClearFormats()
$ws.Range("C115").Value = "Gian Marco Ortenzi"
$ws.Range("D115").Value = "Michele Ribecai"
$ws.Range("E115").Value = "Gana Gian Marco Ortenzi"
$ws.Range("F115").Value = 4.5

# Row 116
$ws.Range("A116").Value = 14600105
$ws.Range("B116").NumberFormat = "@"
$ws.Range("B116").Value = "2025-09-03"
$ws.Range("B116").ClearFormats()
$ws.Range("C116").Value = "Lorenzo Bocchi"
$ws.Range("D116").Value = "Kasra Rahmani"
$ws.Range("E116").Value = "Gana Lorenzo Bocchi"
$ws.Range("F116").Value = 3

# Row 117
$ws.Range("A117").Value = 14600103
$ws.Range("B117").NumberFormat = "@"
$ws.Range("B117").Value = "2025-09-03"
$ws.Range("B117").ClearFormats()
$ws.Range("C117").Value = "Pietro Romeo Scomparin"
$ws.Range("D117").Value = "Leonardo Iemmi"
$ws.Range("E117").Value = "Gana Leonardo Iemmi"
$ws.Range("F117").Value = 5.5

# Row 118
$ws.Range("A118").Value = 14600099
$ws.Range("B118").NumberFormat = "@"
$ws.Range("B118").Value = "2025-09-03"
$ws.Range("B118").ClearFormats()
$ws.Range("C118").Value = "Kirill Kivattsev"
$ws.Range("D118").Value = "Andrea Meduri"
$ws.Range("E118").Value = "Gana Andrea Meduri"
$ws.Range("F118").Value = 3.5

# --- Update sheet dimension to reflect full used range ---
# (Excel COM automatically recalculates dimension based on UsedRange)